$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)  # ALC
$ws.Range("H15").Value = 1479.6451
$ws.Range("I15").Value = 1479.6451
$ws.Range("K15").Value = 4438.9353
$ws.Range("M15").Value = -4269.9353

$ws.Range("H70").Value = 4542.7856

$ws.Range("H73").Value = 4542.7856

$ws.Range("H76").Value = 6305.9375
$ws.Range("I76").Value = 4128.143
$ws.Range("J76").Value = 7999.778
$ws.Range("K76").Value = 4128.143
$ws.Range("L76").Value = 7999.778
$ws.Range("M76").Value = -3813.143
$ws.Range("N76").Value = -8629.778

$ws.Range("H79").Value = 6305.9375
$ws.Range("I79").Value = 4128.143
$ws.Range("J79").Value = 7999.778
$ws.Range("K79").Value = 4128.143
$ws.Range("L79").Value = 7999.778
$ws.Range("M79").Value = -3036.143
$ws.Range("N79").Value = -10183.778

$ws.Range("H100").Value = 5887.846
$ws.Range("I100").Value = 4512.7856
$ws.Range("K100").Value = 4512.7856
$ws.Range("M100").Value = -3971.7856

$ws.Range("H106").Value = 2201.4375
$ws.Range("I106").Value = 2322.3
$ws.Range("K106").Value = 2322.3
$ws.Range("M106").Value = -1691.3

$ws.Range("H116").Value = 4723.7144
$ws.Range("I116").Value = 4558.6
$ws.Range("K116").Value = 4558.6
$ws.Range("M116").Value = -1116.6

$ws.Range("H137").Value = 4331.2617
$ws.Range("I137").Value = 4727.1875
$ws.Range("K137").Value = 14181.5625
$ws.Range("M137").Value = -11631.5625

$ws.Range("H138").Value = 3342.0322
$ws.Range("J138").Value = 4569.184
$ws.Range("L138").Value = 13707.552
$ws.Range("N138").Value = -23987.552

$ws.Range("H141").Value = 3162.3215
$ws.Range("I141").Value = 1663.0435
$ws.Range("J141").Value = 10059
$ws.Range("K141").Value = 4989.1305
$ws.Range("L141").Value = 30177
$ws.Range("M141").Value = 190.8694999999998
$ws.Range("N141").Value = -40537

$ws = $wb.Worksheets.Item(2)  # ARM
$ws.Range("H40").Value = 29747.5
$ws.Range("J40").Value = 29747.5
$ws.Range("L40").Value = 29747.5
$ws.Range("N40").Value = -30099.5

$ws.Range("H61").Value = 1577.4517
$ws.Range("I61").Value = 1562.6666
$ws.Range("K61").Value = 1562.6666
$ws.Range("M61").Value = -1350.6666

$ws.Range("H74").Value = 1621.4584
$ws.Range("I74").Value = 1377.619
$ws.Range("K74").Value = 1377.619
$ws.Range("M74").Value = -503.6189999999999

$ws.Range("H77").Value = 1621.4584
$ws.Range("I77").Value = 1377.619
$ws.Range("K77").Value = 6888.094999999999
$ws.Range("M77").Value = -2520.094999999999

$ws.Range("H97").Value = 3177090.8
$ws.Range("I97").Value = 1605.92
$ws.Range("K97").Value = 1605.92
$ws.Range("M97").Value = -1109.92

$ws.Range("H102").Value = 25644464
$ws.Range("I102").Value = 3697.9
$ws.Range("K102").Value = 3697.9
$ws.Range("M102").Value = -2075.9

$ws.Range("H122").Value = 2673.6667
$ws.Range("I122").Value = 1763.2174
$ws.Range("K122").Value = 5289.6522
$ws.Range("M122").Value = -2839.6522

$ws.Range("H132").Value = 1609.5938
$ws.Range("I132").Value = 1058
$ws.Range("K132").Value = 3174
$ws.Range("M132").Value = -644

$ws.Range("H136").Value = 1577.4517
$ws.Range("I136").Value = 1562.6666
$ws.Range("K136").Value = 4687.9998
$ws.Range("M136").Value = -2137.9998

$ws = $wb.Worksheets.Item(3)  # BSM
$ws.Range("H20").Value = 61032.59
$ws.Range("I20").Value = 1356.2222
$ws.Range("J20").Value = 128168.5
$ws.Range("K20").Value = 1356.2222
$ws.Range("L20").Value = 128168.5
$ws.Range("M20").Value = -1109.2222
$ws.Range("N20").Value = -128662.5

$ws.Range("H86").Value = 2737.7446
$ws.Range("I86").Value = 2114.5264
$ws.Range("J86").Value = 3160.6428
$ws.Range("K86").Value = 2114.5264
$ws.Range("L86").Value = 3160.6428
$ws.Range("M86").Value = -991.5264000000002
$ws.Range("N86").Value = -5406.6428

$ws.Range("H89").Value = 2737.7446
$ws.Range("I89").Value = 2114.5264
$ws.Range("J89").Value = 3160.6428
$ws.Range("K89").Value = 10572.632
$ws.Range("L89").Value = 15803.214
$ws.Range("M89").Value = -4956.632000000001
$ws.Range("N89").Value = -27035.214

$ws.Range("H99").Value = 26129.176
$ws.Range("I99").Value = 26129.176
$ws.Range("K99").Value = 26129.176
$ws.Range("M99").Value = -24631.176

$ws.Range("H141").Value = 80000
$ws.Range("J141").Value = 80000
$ws.Range("L141").Value = 80000
$ws.Range("N141").Value = -90360

$ws = $wb.Worksheets.Item(4)  # CRP
$ws.Range("H31").Value = 2323.8918
$ws.Range("I31").Value = 1895.9412
$ws.Range("K31").Value = 1895.9412
$ws.Range("M31").Value = -1600.9412

$ws.Range("H34").Value = 2323.8918
$ws.Range("I34").Value = 1895.9412
$ws.Range("K34").Value = 1895.9412
$ws.Range("M34").Value = -1693.9412

$ws.Range("H43").Value = 33791
$ws.Range("J43").Value = 33791
$ws.Range("L43").Value = 33791
$ws.Range("N43").Value = -34159

$ws.Range("H101").Value = 33791
$ws.Range("J101").Value = 33791
$ws.Range("L101").Value = 33791
$ws.Range("N101").Value = -40281

$ws.Range("H102").Value = 27499.5
$ws.Range("J102").Value = 27499.5
$ws.Range("L102").Value = 27499.5
$ws.Range("N102").Value = -32367.5

$ws.Range("H103").Value = 14006
$ws.Range("I103").Value = 15341.333
$ws.Range("J103").Value = 10000
$ws.Range("K103").Value = 15341.333
$ws.Range("L103").Value = 10000
$ws.Range("M103").Value = -14169.333
$ws.Range("N103").Value = -12344

$ws.Range("H134").Value = 2202.672
$ws.Range("I134").Value = 2038.9623
$ws.Range("K134").Value = 6116.8869
$ws.Range("M134").Value = -3581.8869

$ws = $wb.Worksheets.Item(5)  # CUL
$ws.Range("H5").Value = 1766
$ws.Range("I5").Value = 1807.75
$ws.Range("J5").Value = 1710.3334
$ws.Range("K5").Value = 5423.25
$ws.Range("L5").Value = 5131.0002
$ws.Range("M5").Value = -5311.25
$ws.Range("N5").Value = -5355.0002

$ws.Range("H135").Value = 1766
$ws.Range("I135").Value = 1807.75
$ws.Range("J135").Value = 1710.3334
$ws.Range("K135").Value = 16269.75
$ws.Range("L135").Value = 15393.0006
$ws.Range("M135").Value = -13734.75
$ws.Range("N135").Value = -20463.0006

$ws = $wb.Worksheets.Item(6)  # GSM
$ws.Range("H80").Value = 22298766
$ws.Range("I80").Value = 223985.2
$ws.Range("K80").Value = 223985.2
$ws.Range("M80").Value = -222987.2

$ws.Range("H83").Value = 22298766
$ws.Range("I83").Value = 223985.2
$ws.Range("K83").Value = 1119926
$ws.Range("M83").Value = -1114934

$ws.Range("H113").Value = 7690.727
$ws.Range("I113").Value = 5999.6665
$ws.Range("K113").Value = 5999.6665
$ws.Range("M113").Value = -3829.6665

$ws.Range("H122").Value = 4267.6665
$ws.Range("I122").Value = 3174.8572
$ws.Range("K122").Value = 9524.571599999999
$ws.Range("M122").Value = -7074.571599999999

$ws.Range("H132").Value = 2389.349
$ws.Range("I132").Value = 2035.3334
$ws.Range("J132").Value = 3893.9167
$ws.Range("K132").Value = 6106.0002
$ws.Range("L132").Value = 11681.7501
$ws.Range("M132").Value = -3576.0002
$ws.Range("N132").Value = -16741.7501

$ws = $wb.Worksheets.Item(7)  # LTW
$ws.Range("H68").Value = 3925.1177
$ws.Range("I68").Value = 3891.25
$ws.Range("K68").Value = 3891.25
$ws.Range("M68").Value = -3142.25

$ws.Range("H71").Value = 3925.1177
$ws.Range("I71").Value = 3891.25
$ws.Range("K71").Value = 19456.25
$ws.Range("M71").Value = -15712.25

$ws.Range("H132").Value = 3237.806
$ws.Range("I132").Value = 2152.843
$ws.Range("K132").Value = 6458.529
$ws.Range("M132").Value = -3928.529

$ws.Range("H136").Value = 5351.0625
$ws.Range("I136").Value = 5351.0625
$ws.Range("K136").Value = 16053.1875
$ws.Range("M136").Value = -13503.1875

$ws = $wb.Worksheets.Item(8)  # WVR
$ws.Range("H96").Value = 34615.53
$ws.Range("I96").Value = 103711.8
$ws.Range("J96").Value = 5825.4165
$ws.Range("K96").Value = 103711.8
$ws.Range("L96").Value = 5825.4165
$ws.Range("M96").Value = -102338.8
$ws.Range("N96").Value = -8571.416499999999

$ws.Range("H122").Value = 2984.3057
$ws.Range("J122").Value = 4123
$ws.Range("L122").Value = 12369
$ws.Range("N122").Value = -17269

$ws.Range("H132").Value = 1301.3793
$ws.Range("I132").Value = 1301.3793
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3904.1379
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1374.1379
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 1472.1875
$ws.Range("I136").Value = 758.3333
$ws.Range("J136").Value = 3613.75
$ws.Range("K136").Value = 2274.9999
$ws.Range("L136").Value = 10841.25
$ws.Range("M136").Value = 275.0001000000002
$ws.Range("N136").Value = -15941.25
